$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab (workbook.xml <sheet name="..."> changed from
# "o554F-HW15.xpc" to "o554F")
$ws.Name = "o554F"

# Append a new row (16) with another averaged-intensity data point produced
# by the new Gaussian-quadrature export, reusing the existing
# "HexGrid-60degTilt5degRes" label already used on row 15
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.338175734288955
$ws.Range("D16").Value = 0.4285253955053088
$ws.Range("E16").Value = 1.047769459359357
$ws.Range("F16").Value = 1.338175734288955
$ws.Range("G16").Value = 0.6975313767768457
$ws.Range("H16").Value = 1.123707654139102
$ws.Range("I16").Value = 1.124106216138495
$ws.Range("J16").Value = 0.4285253955053088
$ws.Range("K16").Value = 0.7381474274323327
$ws.Range("L16").Value = 1.038161580860644
$ws.Range("M16").Value = 0.9599693060346773

# Copy the formatting (border/bold/centered style) of the A column label
# cells down onto the newly-added A16 cell, matching A3:A15
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = 0
